$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (H1) into the two new header cells (I1, J1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Set header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Column I ("I0") is a constant value of 1 for every data row
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("I8").Value = 1

# Column J ("IF") mirrors column H ("IP") values
$ws.Range("J2").Value = $ws.Range("H2").Value2
$ws.Range("J3").Value = $ws.Range("H3").Value2
$ws.Range("J4").Value = $ws.Range("H4").Value2
$ws.Range("J5").Value = $ws.Range("H5").Value2
$ws.Range("J6").Value = $ws.Range("H6").Value2
$ws.Range("J7").Value = $ws.Range("H7").Value2
$ws.Range("J8").Value = $ws.Range("H8").Value2
